# Build site at 2022-09-26 16:07:08 UTC
# Update the LOT2051 "ementa" worksheet:
#  - Row 10 (Objetivos:) B/C content replaced with the professor record that used
#    to live on its own row.
#  - The old stand-alone "101761 - Arnaldo Marcio Ramalho Prata" row (row 13,
#    which had no label in column A) is removed, shifting all following rows
#    up by one.
#  - After the shift, several rows end up showing content that was already
#    present elsewhere in the sheet (a quirk of the source data), so their
#    B/C values are set explicitly to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Objetivos: row now shows the professor info instead of the goals text.
$ws.Range("B10").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C10").Value = "101761 - Arnaldo Márcio Ramalho Prata"

# 2. Remove the old standalone professor-info row (row 13); everything below
#    shifts up by one row.
$ws.Rows.Item(13).Delete()

# 3. Fix up the values on the rows that shifted, per the target content.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2018" already exists verbatim as text in B8/C8 - copy it across so it
# isn't re-parsed into a date serial number by a fresh literal assignment.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("B18").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C18").Value = "101761 - Arnaldo Márcio Ramalho Prata"

$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = [P1 +(2 x P2)] / 3"
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = [P1 +(2 x P2)] / 3"

$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
